$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 612.7308
$ws.Range("I19").Value = 123.818184
$ws.Range("K19").Value = 123.818184
$ws.Range("M19").Value = 51.181816
$ws.Range("H62").Value = 7231.0835
$ws.Range("I62").Value = 6681.857
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 6681.857
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -6057.857
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 7231.0835
$ws.Range("I65").Value = 6681.857
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 33409.285
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -30289.285
$ws.Range("N65").Value = -46240
$ws.Range("H86").Value = 6409.579
$ws.Range("I86").Value = 4768.1113
$ws.Range("K86").Value = 4768.1113
$ws.Range("M86").Value = -3645.1113
$ws.Range("H89").Value = 6409.579
$ws.Range("I89").Value = 4768.1113
$ws.Range("K89").Value = 23840.5565
$ws.Range("M89").Value = -18224.5565
$ws.Range("H100").Value = 4099.8887
$ws.Range("I100").Value = 4284.8335
$ws.Range("K100").Value = 4284.8335
$ws.Range("M100").Value = -3743.8335
$ws.Range("H113").Value = 7151.2
$ws.Range("I113").Value = 8490.799999999999
$ws.Range("J113").Value = 6481.4
$ws.Range("K113").Value = 8490.799999999999
$ws.Range("L113").Value = 6481.4
$ws.Range("M113").Value = -5236.799999999999
$ws.Range("N113").Value = -12989.4
$ws.Range("H116").Value = 6514.75
$ws.Range("I116").Value = 4725.727
$ws.Range("J116").Value = 8701.333000000001
$ws.Range("K116").Value = 4725.727
$ws.Range("L116").Value = 8701.333000000001
$ws.Range("M116").Value = -1283.727
$ws.Range("N116").Value = -15585.333
$ws.Range("H132").Value = 27029144
$ws.Range("I132").Value = 28573614
$ws.Range("J132").Value = 953
$ws.Range("K132").Value = 85720842
$ws.Range("L132").Value = 2859
$ws.Range("M132").Value = -85718312
$ws.Range("N132").Value = -7919

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2095911.9
$ws.Range("I2").Value = 2977762.5
$ws.Range("J2").Value = 1516.625
$ws.Range("K2").Value = 2977762.5
$ws.Range("L2").Value = 1516.625
$ws.Range("M2").Value = -2977649.5
$ws.Range("N2").Value = -1742.625
$ws.Range("H32").Value = 8713.41
$ws.Range("I32").Value = 4880.7964
$ws.Range("K32").Value = 4880.7964
$ws.Range("M32").Value = -4593.7964
$ws.Range("H45").Value = 8551858
$ws.Range("I45").Value = 13988770
$ws.Range("K45").Value = 13988770
$ws.Range("M45").Value = -13988393
$ws.Range("H61").Value = 3456.7585
$ws.Range("I61").Value = 3343
$ws.Range("J61").Value = 4992.5
$ws.Range("K61").Value = 3343
$ws.Range("L61").Value = 4992.5
$ws.Range("M61").Value = -3131
$ws.Range("N61").Value = -5416.5
$ws.Range("H63").Value = 4602.25
$ws.Range("I63").Value = 2313.7273
$ws.Range("J63").Value = 7399.3335
$ws.Range("K63").Value = 2313.7273
$ws.Range("L63").Value = 7399.3335
$ws.Range("M63").Value = -1627.7273
$ws.Range("N63").Value = -8771.333500000001
$ws.Range("H66").Value = 4602.25
$ws.Range("I66").Value = 2313.7273
$ws.Range("J66").Value = 7399.3335
$ws.Range("K66").Value = 11568.6365
$ws.Range("L66").Value = 36996.6675
$ws.Range("M66").Value = -8136.636500000001
$ws.Range("N66").Value = -43860.6675
$ws.Range("H102").Value = 2454398.2
$ws.Range("I102").Value = 2979279
$ws.Range("J102").Value = 4955
$ws.Range("K102").Value = 2979279
$ws.Range("L102").Value = 4955
$ws.Range("M102").Value = -2977657
$ws.Range("N102").Value = -8199
$ws.Range("H110").Value = 695453.9399999999
$ws.Range("I110").Value = 732019.9399999999
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 732019.9399999999
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = -729974.9399999999
$ws.Range("N110").Value = -4790
$ws.Range("H116").Value = 2095911.9
$ws.Range("I116").Value = 2977762.5
$ws.Range("J116").Value = 1516.625
$ws.Range("K116").Value = 2977762.5
$ws.Range("L116").Value = 1516.625
$ws.Range("M116").Value = -2975468.5
$ws.Range("N116").Value = -6104.625
$ws.Range("H132").Value = 2465.8948
$ws.Range("I132").Value = 1843.8235
$ws.Range("J132").Value = 7753.5
$ws.Range("K132").Value = 5531.470499999999
$ws.Range("L132").Value = 23260.5
$ws.Range("M132").Value = -3001.470499999999
$ws.Range("N132").Value = -28320.5
$ws.Range("H136").Value = 3456.7585
$ws.Range("I136").Value = 3343
$ws.Range("J136").Value = 4992.5
$ws.Range("K136").Value = 10029
$ws.Range("L136").Value = 14977.5
$ws.Range("M136").Value = -7479
$ws.Range("N136").Value = -20077.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2095911.9
$ws.Range("I3").Value = 2977762.5
$ws.Range("J3").Value = 1516.625
$ws.Range("K3").Value = 2977762.5
$ws.Range("L3").Value = 1516.625
$ws.Range("M3").Value = -2977648.5
$ws.Range("N3").Value = -1744.625
$ws.Range("H38").Value = 88571.28999999999
$ws.Range("J38").Value = 88571.28999999999
$ws.Range("L38").Value = 88571.28999999999
$ws.Range("N38").Value = -89403.28999999999
$ws.Range("H80").Value = 389.03125
$ws.Range("I80").Value = 324.45456
$ws.Range("J80").Value = 422.85715
$ws.Range("K80").Value = 324.45456
$ws.Range("L80").Value = 422.85715
$ws.Range("M80").Value = 673.54544
$ws.Range("N80").Value = -2418.85715
$ws.Range("H83").Value = 389.03125
$ws.Range("I83").Value = 324.45456
$ws.Range("J83").Value = 422.85715
$ws.Range("K83").Value = 1622.2728
$ws.Range("L83").Value = 2114.28575
$ws.Range("M83").Value = 3369.7272
$ws.Range("N83").Value = -12098.28575
$ws.Range("H94").Value = 6263124
$ws.Range("I94").Value = 11111837
$ws.Range("K94").Value = 11111837
$ws.Range("M94").Value = -11111386
$ws.Range("H99").Value = 8932920
$ws.Range("I99").Value = 11908769
$ws.Range("K99").Value = 11908769
$ws.Range("M99").Value = -11907271
$ws.Range("H105").Value = 5682683.5
$ws.Range("I105").Value = 6945300
$ws.Range("K105").Value = 6945300
$ws.Range("M105").Value = -6943553
$ws.Range("H107").Value = 3762185.2
$ws.Range("I107").Value = 4763968
$ws.Range("K107").Value = 4763968
$ws.Range("M107").Value = -4762048
$ws.Range("H134").Value = 4392.8374
$ws.Range("I134").Value = 3350.2727
$ws.Range("K134").Value = 10050.8181
$ws.Range("M134").Value = -7515.8181

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1170.125
$ws.Range("I16").Value = 924.0833
$ws.Range("K16").Value = 924.0833
$ws.Range("M16").Value = -637.0833
$ws.Range("H107").Value = 1472.421
$ws.Range("I107").Value = 810.375
$ws.Range("K107").Value = 810.375
$ws.Range("M107").Value = 1109.625
$ws.Range("H113").Value = 1170.125
$ws.Range("I113").Value = 924.0833
$ws.Range("K113").Value = 924.0833
$ws.Range("M113").Value = 1245.9167

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 53531.105
$ws.Range("I5").Value = 807.3077
$ws.Range("K5").Value = 2421.9231
$ws.Range("M5").Value = -2309.9231
$ws.Range("H32").Value = 900000000
$ws.Range("I32").Value = 1000000000
$ws.Range("K32").Value = 3000000000
$ws.Range("M32").Value = -2999999717
$ws.Range("H46").Value = 240955.42
$ws.Range("I46").Value = 1666791.5
$ws.Range("K46").Value = 5000374.5
$ws.Range("M46").Value = -5000283.5
$ws.Range("H107").Value = 1975.1111
$ws.Range("I107").Value = 2694.75
$ws.Range("J107").Value = 1399.4
$ws.Range("K107").Value = 8084.25
$ws.Range("L107").Value = 4198.200000000001
$ws.Range("M107").Value = -6164.25
$ws.Range("N107").Value = -8038.200000000001
$ws.Range("H135").Value = 53531.105
$ws.Range("I135").Value = 807.3077
$ws.Range("K135").Value = 7265.7693
$ws.Range("M135").Value = -4730.7693

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 882644.75
$ws.Range("I97").Value = 1401251.4
$ws.Range("J97").Value = 1013.4
$ws.Range("K97").Value = 1401251.4
$ws.Range("L97").Value = 1013.4
$ws.Range("M97").Value = -1400755.4
$ws.Range("N97").Value = -2005.4
$ws.Range("H132").Value = 2545.4307
$ws.Range("I132").Value = 2365.3333
$ws.Range("J132").Value = 3201.5
$ws.Range("K132").Value = 7095.999899999999
$ws.Range("L132").Value = 9604.5
$ws.Range("M132").Value = -4565.999899999999
$ws.Range("N132").Value = -14664.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 24998
$ws.Range("I29").Value = 24998
$ws.Range("K29").Value = 24998
$ws.Range("M29").Value = -24703
$ws.Range("H46").Value = 5375.846
$ws.Range("I46").Value = 2277.6
$ws.Range("K46").Value = 2277.6
$ws.Range("M46").Value = -2089.6
$ws.Range("H53").Value = 30143.666
$ws.Range("I53").Value = 6329
$ws.Range("J53").Value = 42051
$ws.Range("K53").Value = 6329
$ws.Range("L53").Value = 42051
$ws.Range("M53").Value = -5811
$ws.Range("N53").Value = -43087

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M37").ClearContents()
$ws.Range("H31").Value = 16803.8
$ws.Range("I31").Value = 10000
$ws.Range("K31").Value = 10000
$ws.Range("M31").Value = -9652
$ws.Range("H37").Value = 28000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 28000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 28000
$ws.Range("N37").Value = -28406
$ws.Range("H42").Value = 36998
$ws.Range("I42").Value = 36998
$ws.Range("K42").Value = 36998
$ws.Range("M42").Value = -36620
$ws.Range("H107").Value = 66668144
$ws.Range("I107").Value = 83333780
$ws.Range("J107").Value = 5617.6665
$ws.Range("K107").Value = 250001340
$ws.Range("L107").Value = 16852.9995
$ws.Range("M107").Value = -249999420
$ws.Range("N107").Value = -20692.9995
